# Insert a new weekly record for Red Globe / Segunda (Región de O'Higgins)
# at row 32, shifting every subsequent row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 32, pushing rows 32..90 to 33..91.
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new record's data.
$ws.Cells.Item(32, 1).Value = 11
$ws.Cells.Item(32, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(32, 3).Value = "Bíobío"
$ws.Cells.Item(32, 4).Value = 44421
$ws.Cells.Item(32, 5).Value = 8
$ws.Cells.Item(32, 6).Value = "Fruta"
$ws.Cells.Item(32, 7).Value = 100109
$ws.Cells.Item(32, 8).Value = "Uva"
$ws.Cells.Item(32, 9).Value = 100109001
$ws.Cells.Item(32, 10).Value = "Uva"
$ws.Cells.Item(32, 11).Value = "Red Globe"
$ws.Cells.Item(32, 12).Value = "Segunda"
$ws.Cells.Item(32, 13).Value = 100
$ws.Cells.Item(32, 14).Value = 10000
$ws.Cells.Item(32, 15).Value = 11000
$ws.Cells.Item(32, 16).Value = 10500
$ws.Cells.Item(32, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(32, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(32, 19).Value = 1312
$ws.Cells.Item(32, 20).Value = 8
